$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E24").Value = 8
$ws.Range("E25").Value = 8.99
$ws.Range("E26").Value = 9.640000000000001
try { $wb.RefreshAll() ; Write-Host "RefreshAll ok" } catch { Write-Host "RefreshAll err: $_" }
